$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.658.52'
$ws.Range("E2").Value = '  +2.09%  '

$ws.Range("D3").Value = '3.945.87'
$ws.Range("E3").Value = '  +0.66%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '517.72'
$ws.Range("E5").Value = '  +6.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.93'
$ws.Range("E6").Value = '  -1.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("E7").Value = '  -0.22%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.732'
$ws.Range("E9").Value = '  -0.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.173'
$ws.Range("E10").Value = '  +4.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000345'
$ws.Range("E11").Value = '  -1.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.25'
$ws.Range("E12").Value = '  +0.46%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.43'
$ws.Range("E13").Value = '  -3.66%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '4.560.22'
$ws.Range("E14").Value = '  +0.29%  '

$ws.Range("D15").Value = '3.937.82'
$ws.Range("E15").Value = '  +0.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.21'
$ws.Range("E16").Value = '  -1.37%  '

$ws.Range("E17").Value = '  -0.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.23'
$ws.Range("E18").Value = '  +7.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '19.88'
$ws.Range("E19").Value = '  -0.28%  '

$ws.Range("D20").Value = '69.529.37'
$ws.Range("E20").Value = '  +1.66%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '434.53'
$ws.Range("E21").Value = '  -1.86%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.42'
$ws.Range("E22").Value = '  -1.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.58'
$ws.Range("E23").Value = '  -4.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.41'
$ws.Range("E24").Value = '  +0.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.87'
$ws.Range("E25").Value = '  +3.21%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.91'
$ws.Range("E26").Value = '  +7.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.17'
$ws.Range("E27").Value = '  -3.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '36.95'
$ws.Range("E28").Value = '  -4.32%  '

$ws.Range("E29").Value = '  -1.44%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '702.27'
$ws.Range("E30").Value = '  -2.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.33'
$ws.Range("E31").Value = '  -3.68%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.128'
$ws.Range("E32").Value = '  -2.62%  '

$ws.Range("E33").Value = '  -1.52%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '67.68'
$ws.Range("E34").Value = '  +10.42%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.443'
$ws.Range("E35").Value = '  +7.64%  '

$ws.Range("D36").Value = '0.0₃0884'
$ws.Range("E36").Value = '  +1.17%  '

$ws.Range("E37").Value = '  -5.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '40.43'
$ws.Range("E38").Value = '  -4.35%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.149'
$ws.Range("E39").Value = '  -0.12%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0486'
$ws.Range("E42").Value = '  +0.87%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.83'
$ws.Range("E43").Value = '  -6.46%  '

$ws.Range("E44").Value = '  +6.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.03'
$ws.Range("E45").Value = '  -7.17%  '

$ws.Range("E46").Value = '  +0.96%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.37'
$ws.Range("E47").Value = '  +1.89%  '

$ws.Range("D48").Value = '0.0₆0357'
$ws.Range("E48").Value = '  +2.35%  '

$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.97'
$ws.Range("E49").Value = '  +4.80%  '

$ws.Range("B50").Value = 'LidoDAOToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.36'
$ws.Range("E50").Value = '  -2.07%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.10'
$ws.Range("E51").Value = '  -2.26%  '
